{"js": "// New math-expression text for each table cell, in row-major document\n// order (row 0 cols 0..4, row 1 cols 0..4, ...). This corresponds 1:1 to\n// the <w:t> replacements in the target diff.\nconst newValues = [\"60+20=\", \"80+14=\", \"20+23=\", \"33+66=\", \"29+40=\", \"51+28=\", \"43+21=\", \"55-11=\", \"3+18=\", \"20+16=\", \"5+12=\", \"43-5=\", \"81+14=\", \"77-7=\", \"99-17=\", \"61-25=\", \"46-40=\", \"13+83=\", \"9+37=\", \"78-6=\", \"91-7=\", \"88+3=\", \"34+41=\", \"29+36=\", \"87-8=\", \"81-25=\", \"64+20=\", \"93-70=\", \"40-12=\", \"10+2=\", \"22+77=\", \"97-90=\", \"42-26=\", \"39-26=\", \"27+69=\", \"19+7=\", \"81-14=\", \"74-4=\", \"13-7=\", \"47+13=\", \"13+11=\", \"9+1=\", \"93-60=\", \"49+21=\", \"83-5=\", \"60+19=\", \"31+16=\", \"44-9=\", \"66+26=\", \"7-3=\", \"17+19=\", \"26-2=\", \"79-16=\", \"42+55=\", \"65-35=\", \"98-69=\", \"89-89=\", \"32+54=\", \"39+20=\", \"71-51=\", \"30-16=\", \"42+35=\", \"81-7=\", \"93-61=\", \"25+62=\", \"86-49=\", \"88-22=\", \"51-15=\", \"54+31=\", \"59+4=\", \"87-6=\", \"93-43=\", \"20+68=\", \"99-65=\", \"81-1=\", \"36-13=\", \"70-38=\", \"92-0=\", \"80-35=\", \"85+11=\", \"33-27=\", \"51-32=\", \"16+38=\", \"77-10=\", \"64-13=\", \"38+61=\", \"22+51=\", \"28+15=\", \"73+0=\", \"32+61=\", \"14+28=\", \"63-19=\", \"53+5=\", \"10+34=\", \"87-29=\", \"75-61=\", \"64-53=\", \"70-23=\", \"33+0=\", \"24+36=\"];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    if (idx >= newValues.length) break;\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# New math-expression text for each table cell, in row-major document\n# order (row 0 cols 0..4, row 1 cols 0..4, ...). This corresponds 1:1 to\n# the <w:t> replacements in the target diff.\n$newValues = @(\n    @(\"60+20=\", \"80+14=\", \"20+23=\", \"33+66=\", \"29+40=\"),\n    @(\"51+28=\", \"43+21=\", \"55-11=\", \"3+18=\", \"20+16=\"),\n    @(\"5+12=\", \"43-5=\", \"81+14=\", \"77-7=\", \"99-17=\"),\n    @(\"61-25=\", \"46-40=\", \"13+83=\", \"9+37=\", \"78-6=\"),\n    @(\"91-7=\", \"88+3=\", \"34+41=\", \"29+36=\", \"87-8=\"),\n    @(\"81-25=\", \"64+20=\", \"93-70=\", \"40-12=\", \"10+2=\"),\n    @(\"22+77=\", \"97-90=\", \"42-26=\", \"39-26=\", \"27+69=\"),\n    @(\"19+7=\", \"81-14=\", \"74-4=\", \"13-7=\", \"47+13=\"),\n    @(\"13+11=\", \"9+1=\", \"93-60=\", \"49+21=\", \"83-5=\"),\n    @(\"60+19=\", \"31+16=\", \"44-9=\", \"66+26=\", \"7-3=\"),\n    @(\"17+19=\", \"26-2=\", \"79-16=\", \"42+55=\", \"65-35=\"),\n    @(\"98-69=\", \"89-89=\", \"32+54=\", \"39+20=\", \"71-51=\"),\n    @(\"30-16=\", \"42+35=\", \"81-7=\", \"93-61=\", \"25+62=\"),\n    @(\"86-49=\", \"88-22=\", \"51-15=\", \"54+31=\", \"59+4=\"),\n    @(\"87-6=\", \"93-43=\", \"20+68=\", \"99-65=\", \"81-1=\"),\n    @(\"36-13=\", \"70-38=\", \"92-0=\", \"80-35=\", \"85+11=\"),\n    @(\"33-27=\", \"51-32=\", \"16+38=\", \"77-10=\", \"64-13=\"),\n    @(\"38+61=\", \"22+51=\", \"28+15=\", \"73+0=\", \"32+61=\"),\n    @(\"14+28=\", \"63-19=\", \"53+5=\", \"10+34=\", \"87-29=\"),\n    @(\"75-61=\", \"64-53=\", \"70-23=\", \"33+0=\", \"24+36=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r-1][$c-1]\n    }\n}\n"}
